# Update stats for 2026-02 (row 27 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6565
$ws.Range("D27").Value = 6131174
$ws.Range("E27").Value = 933.9183549124143
$ws.Range("F27").Value = 10.33613445378152
$ws.Range("H27").Value = 25.75459331716681
